$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 233, shifting existing rows 233-272 down to 234-273
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row with the new data record
$ws.Cells.Item(233, 1).Value = 10
$ws.Cells.Item(233, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(233, 3).Value = "La Araucanía"
$ws.Cells.Item(233, 4).Value = 44511
$ws.Cells.Item(233, 5).Value = 9
$ws.Cells.Item(233, 6).Value = 100112024
$ws.Cells.Item(233, 7).Value = "Choclo"
$ws.Cells.Item(233, 8).Value = "Dulce o Americano"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 100
$ws.Cells.Item(233, 11).Value = 28000
$ws.Cells.Item(233, 12).Value = 28000
$ws.Cells.Item(233, 13).Value = 28000
$ws.Cells.Item(233, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(233, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(233, 16).Value = 400
$ws.Cells.Item(233, 17).Value = 70
$ws.Cells.Item(233, 18).Value = "Hortaliza"

# Match the date-formatted style used by column D elsewhere in the table
$ws.Cells.Item(233, 4).NumberFormat = $ws.Cells.Item(234, 4).NumberFormat
